$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: copy the existing header formatting (bold, bordered,
# centered) from AC1 onto AD1:AF1, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record columns for every player row: Wins=68, Losses=94, Ties=0.
for ($r = 2; $r -le 48; $r++) {
    $ws.Range("AD$r").Value = 68
    $ws.Range("AE$r").Value = 94
    $ws.Range("AF$r").Value = 0
}
